$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.219.89"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.784.51"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3433"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.33"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07476"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.04%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.71"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.460"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "1.786.29"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.077"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06667"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.650"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("D23").Value = "27.227.42"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.39"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.420"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.504"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.540"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.38"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "1.988.27"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.96"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.011"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.081"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08673"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.30%  "
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6929"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.465"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06338"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2198"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.790"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02339"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.242"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.41"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6492"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.856"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.142"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.76"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07123"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.25%  "
